# Update "想去人数" (number of people interested) counts for several rows
# on both the "展览" and "全部类型" worksheets.
#
# Row -> new value mapping (column F):
#   F3: 71   -> 74
#   F4: 481  -> 482
#   F5: 4676 -> 4685
#   F6: 365  -> 366
#   F8: 287  -> 288
#   F9: 724  -> 726

$wb = $excel.ActiveWorkbook

$updates = @{
    3 = 74
    4 = 482
    5 = 4685
    6 = 366
    8 = 288
    9 = 726
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
